$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# like "187.70" or "0.550" are stored verbatim instead of being
# coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '74.866.38'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '2.821.68'
$ws.Range("E3").Value = '  +7.29%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '187.70'
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("D6").Value = '595.04'
$ws.Range("E6").Value = '  +2.11%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '0.550'
$ws.Range("E8").Value = '  +3.03%  '
$ws.Range("D9").Value = '0.193'
$ws.Range("E9").Value = '  -4.93%  '
$ws.Range("D10").Value = '2.820.81'
$ws.Range("E10").Value = '  +7.32%  '
$ws.Range("D11").Value = '0.160'
$ws.Range("E11").Value = '  -1.32%  '
$ws.Range("E12").Value = '  +3.61%  '
$ws.Range("E13").Value = '  +2.89%  '
$ws.Range("D14").Value = '3.342.14'
$ws.Range("E14").Value = '  +7.50%  '
$ws.Range("D15").Value = '74.809.86'
$ws.Range("E15").Value = '  +1.10%  '
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("D17").Value = '26.82'
$ws.Range("E17").Value = '  +2.14%  '
$ws.Range("D18").Value = '2.820.68'
$ws.Range("E18").Value = '  +7.21%  '
$ws.Range("E19").Value = '  -1.87%  '
$ws.Range("D20").Value = '12.31'
$ws.Range("E20").Value = '  +3.99%  '
$ws.Range("D21").Value = '377.34'
$ws.Range("E21").Value = '  +1.35%  '
$ws.Range("D22").Value = '2.27'
$ws.Range("E22").Value = '  -1.74%  '
$ws.Range("E23").Value = '  -0.73%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  +1.27%  '
$ws.Range("E26").Value = '  +7.37%  '
$ws.Range("D27").Value = '4.17'
$ws.Range("E27").Value = '  +0.48%  '
$ws.Range("D28").Value = '9.73'
$ws.Range("E28").Value = '  +3.60%  '
$ws.Range("E29").Value = '  +10.11%  '
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").Value = '517.74'
$ws.Range("E31").Value = '  -1.95%  '
$ws.Range("E32").Value = '  -0.53%  '
$ws.Range("E33").Value = '  +0.64%  '
$ws.Range("D34").Value = '1.78'
$ws.Range("E34").Value = '  +2.55%  '
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").Value = '163.47'
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").Value = '19.92'
$ws.Range("E37").Value = '  +3.89%  '
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("D39").Value = '19.36'
$ws.Range("E39").Value = '  +0.53%  '
$ws.Range("D40").Value = '186.93'
$ws.Range("E40").Value = '  +16.15%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("E42").Value = '  +3.39%  '
$ws.Range("E43").Value = '  +1.49%  '
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("E45").Value = '  +1.59%  '
$ws.Range("E46").Value = '  +2.66%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.0856'
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").Value = '2.32'
$ws.Range("E48").Value = '  -2.14%  '
$ws.Range("E49").Value = '  +9.49%  '
$ws.Range("E50").Value = '  +2.46%  '
$ws.Range("D51").Value = '0.636'
$ws.Range("E51").Value = '  +8.26%  '

# Restore the original (default) cell style for column D now that
# the text values are committed, so no stray number-format style
# lingers on cells that did not have one before.
$ws.Range("D2:D51").Style = "Normal"
